$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to Text format before writing so Excel does not
# silently reinterpret numeric-looking / percentage-looking strings as
# numbers, then restore the default ("Normal") style so no stray
# number-format styling is left behind on the cell.
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '323.83'
Set-TextValue 'E2' '-1.78%'
Set-TextValue 'D3' '39.60'
Set-TextValue 'E3' '-0.91%'
Set-TextValue 'D4' '5.858'
Set-TextValue 'E4' '11.13%'
Set-TextValue 'D5' '0.08040'
Set-TextValue 'E5' '-0.65%'
Set-TextValue 'D6' '2.019'
Set-TextValue 'E6' '4.65%'
Set-TextValue 'B7' 'KuCoinToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D7' '8.633'
Set-TextValue 'E7' '-0.16%'
Set-TextValue 'B8' 'BTSEToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D8' '2.948'
Set-TextValue 'E8' '-0.33%'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9307'
Set-TextValue 'E9' '-0.70%'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1268'
Set-TextValue 'E10' '-6.50%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1954'
Set-TextValue 'E11' '-1.27%'
Set-TextValue 'B12' 'MCDex'
Set-TextValue 'C12' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D12' '8.698'
Set-TextValue 'E12' '28.52%'
Set-TextValue 'B13' 'MandalaExchangeToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D13' '0.09147'
Set-TextValue 'E13' '0.63%'
Set-TextValue 'B14' 'BitrueCoin'
Set-TextValue 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D14' '0.03586'
Set-TextValue 'E14' '2.35%'
Set-TextValue 'B15' 'BitMartToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D15' '0.1046'
Set-TextValue 'E15' '9.12%'
Set-TextValue 'B16' 'BitForexToken'
Set-TextValue 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001292'
Set-TextValue 'E16' '-7.13%'
Set-TextValue 'B17' 'TigerCash'
Set-TextValue 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D17' '0.006176'
Set-TextValue 'E17' '-5.35%'
Set-TextValue 'B18' 'LEO'
Set-TextValue 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D18' '3.352'
Set-TextValue 'E18' '-0.42%'
Set-TextValue 'B19' 'GateToken'
Set-TextValue 'C19' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D19' '4.571'
Set-TextValue 'E19' '1.37%'
Set-TextValue 'E20' '0.47%'
Set-TextValue 'D21' '0.1371'
Set-TextValue 'E21' '4.42%'
Set-TextValue 'D22' '0.2448'
Set-TextValue 'E22' '-4.72%'
Set-TextValue 'D23' '0.04405'
Set-TextValue 'E23' '-0.53%'
Set-TextValue 'D24' '0.001261'
Set-TextValue 'E24' '3.27%'
Set-TextValue 'D25' '0.004388'
Set-TextValue 'E25' '1.92%'
Set-TextValue 'D26' '0.0001150'
Set-TextValue 'E26' '-10.89%'
Set-TextValue 'D39' '0.02541'
Set-TextValue 'E39' '1.88%'
Set-TextValue 'D40' '0.05244'
Set-TextValue 'E40' '0.40%'
Set-TextValue 'D41' '0.007472'
Set-TextValue 'E41' '-3.27%'
Set-TextValue 'D42' '0.009595'
Set-TextValue 'E42' '4.35%'
Set-TextValue 'D43' '0.1407'
Set-TextValue 'E43' '-1.47%'
Set-TextValue 'D44' '0.002116'
Set-TextValue 'E44' '-2.52%'
Set-TextValue 'D45' '0.009980'
Set-TextValue 'E45' '7.09%'
Set-TextValue 'D46' '0.00006745'
Set-TextValue 'E46' '1.45%'
Set-TextValue 'D47' '0.00000000750'
Set-TextValue 'E47' '-0.05%'
Set-TextValue 'D48' '0.003000'
Set-TextValue 'E48' '-9.92%'
Set-TextValue 'D49' '0.002290'
Set-TextValue 'E49' '-7.78%'
Set-TextValue 'D50' '0.00002100'
Set-TextValue 'E50' '-0.05%'
Set-TextValue 'D51' '0.0002000'
Set-TextValue 'E51' '-0.05%'
